# Update objective (B), gap (C) and solve time (D) values for rows 2-11
# on the active worksheet to reflect corrected fixed-recourse data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ B = -71.65924151606299;  C = 0.0005680745486136919; D = 25.557223755  }
    3  = @{ B = -69.45256735952894;  C = 0.0;                   D = 53.479246961  }
    4  = @{ B = -71.24685956643927;  C = 0.05095443095246155;   D = 8.61043671    }
    5  = @{ B = -71.65440618645333;  C = 0.0942196409173878;    D = 5.525219539   }
    6  = @{ B = -70.22823688201865;  C = 0.05353691840979169;   D = 7.38619507    }
    7  = @{ B = -71.84370853874832;  C = 0.0;                   D = 129.531268755 }
    8  = @{ B = -67.12860480761344;  C = 0.000000550528395181048; D = 17.921225371  }
    9  = @{ B = -71.56787798834289;  C = 0.0958882217621773;    D = 4.964546475   }
    10 = @{ B = -70.884938041055;    C = 0.05531266370132536;   D = 6.314387318   }
    11 = @{ B = -68.19423200083996;  C = 0.0;                   D = 36.246093994  }
}

foreach ($row in $values.Keys) {
    $rowData = $values[$row]
    $ws.Range("B$row").Value = $rowData.B
    $ws.Range("C$row").Value = $rowData.C
    $ws.Range("D$row").Value = $rowData.D
}
